$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("transitional_list")

$ws1.Range("B3").Value = "V13I,I17L,T19A,V30I,Y40C,N43H,N47T,I62T,E64K,K64N,I77K,K93R,D147N,V149F,I165V,K187R,E199N,K199E, M241V,R249I,I263V,L269I,V313A"
$ws1.Range("B4").Value = "L22F,N41D,V50A,M51I,K75R,D86S,K93N,D127G,I149V,S161N,I176M,I194V,D199N,F205L,I257V,V263I,I263V,R264H,K267T,S284F,D309N,H310Y,V312I"
$ws1.Range("B5").Value = "I26V,V26I,D41N,N43D,Q49H,L52F,K62T,I73V,L81P,A82T,I257V,I263V,I312T"

$ws1.Rows.Item(3).RowHeight = 47.25

$ws1.Range("B3").WrapText = $true
$ws1.Range("B3").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop

$ws1.Range("A1:B5").Select()
